# Apply odds updates to Sheet1 for the 2025-12-23 Betfair Back/Lay workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell address -> new numeric value, derived from the OOXML diff.
$updates = @{
    "F2"  = 1.68
    "G2"  = 1.69
    "T2"  = 1.86
    "W2"  = 2.44

    "F3"  = 1.39
    "G3"  = 1.49
    "H3"  = 9.6
    "J3"  = 3.85
    "K3"  = 5.3
    "L3"  = 1.4
    "O3"  = 1.44
    "S3"  = 4.5
    "W3"  = 3
    "X3"  = 11.5
    "AC3" = 13
    "AD3" = 60
    "AF3" = 7.8
    "AG3" = 13.5
    "AK3" = 25
    "AN3" = 11

    "H4"  = 1.41
    "I4"  = 1.43
    "J4"  = 5.2
    "K4"  = 5.7
    "T4"  = 2.12
    "U4"  = 1.78
    "V4"  = 3.25
    "X4"  = 18.5
    "Y4"  = 8.4
    "Z4"  = 7.8
    "AE4" = 15.5
    "AF4" = 85
    "AH4" = 29
    "AI4" = 42
    "AJ4" = 400
    "AK4" = 180
    "AL4" = 150
    "AM4" = 210
    "AN4" = 250
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$wb.Save()
